$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 23
$ws.Range("E3").Value = 18
$ws.Range("F6").Value = 24
$ws.Range("H6").Value = 24
$ws.Range("F10").Value = 8
$ws.Range("H10").Value = 8
$ws.Range("F12").Value = 7
$ws.Range("H12").Value = 7
$ws.Range("E16").Value = 283
$ws.Range("E18").Value = 86
